$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$row = 54
$ws.Cells.Item($row, 1).Value = 1759778989
$ws.Cells.Item($row, 2).Value = "update"
$ws.Cells.Item($row, 3).Value = "doc"
$ws.Cells.Item($row, 4).Value = "bevnat-info"
$ws.Cells.Item($row, 6).Value = "last_update"

# G54 and H54 hold digit-only strings ("1706219962" / "1706239962") that must be
# stored as text (shared string), not auto-coerced into numbers, and without
# mutating the cell's style (no NumberFormat "@"/quote-prefix). Using a helper
# cell with a TEXT() formula and pasting only the value achieves a true text
# value while leaving the destination cell's default style untouched.
$helper = $ws.Cells.Item(1000, 1)

$helper.Formula = '=TEXT(1706219962,"0")'
$helper.Copy() | Out-Null
$ws.Cells.Item($row, 7).PasteSpecial(-4163) | Out-Null

$helper.Formula = '=TEXT(1706239962,"0")'
$helper.Copy() | Out-Null
$ws.Cells.Item($row, 8).PasteSpecial(-4163) | Out-Null

$helper.ClearContents() | Out-Null
$excel.CutCopyMode = 0
